# Daily COVID stats update: st 01. 07. 2021
# 1) Row 89 previously had stray AgTests/AgPosit (F/G) values that don't belong
#    to that date (the Ag-testing series starts later) - remove them.
# 2) Refresh the cumulative AgTests (F) / AgPosit (G) columns for the trailing
#    block of existing rows with corrected totals.
# 3) Append the new day's row (44377 = 2021-06-30) with the day's figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Clear the erroneous F89/G89 values ---
$ws.Range("F89:G89").ClearContents()

# --- 2) Update AgTests (F) / AgPosit (G) for the affected rows ---
$updates = @(
    @{Row=272; F=30377; G=1604},
    @{Row=278; F=31353; G=2129},
    @{Row=279; F=42687; G=2982},
    @{Row=280; F=34966; G=2310},
    @{Row=281; F=47131; G=3179},
    @{Row=282; F=46463; G=2737},
    @{Row=285; F=42877; G=3411},
    @{Row=286; F=55175; G=4275},
    @{Row=287; F=58906; G=3710},
    @{Row=288; F=59445; G=3978},
    @{Row=289; F=62981; G=3590},
    @{Row=291; F=15296; G=496},
    @{Row=292; F=83201; G=7278},
    @{Row=293; F=83320; G=5765},
    @{Row=294; F=94238; G=4943},
    @{Row=295; F=17253; G=1043},
    @{Row=296; F=2463; G=141},
    @{Row=299; F=66381; G=6869},
    @{Row=300; F=72864; G=6963},
    @{Row=301; F=72430; G=5697},
    @{Row=302; F=78588; G=5608},
    @{Row=338; F=221315; G=3047},
    @{Row=346; F=675004; G=4822},
    @{Row=393; F=308046; G=1240},
    @{Row=423; F=439590; G=636},
    @{Row=425; F=137521; G=548},
    @{Row=426; F=106736; G=380},
    @{Row=427; F=89027; G=366},
    @{Row=428; F=102355; G=389},
    @{Row=429; F=171220; G=430},
    @{Row=430; F=169419; G=257},
    @{Row=432; F=118194; G=414},
    @{Row=433; F=85970; G=264},
    @{Row=434; F=79049; G=279},
    @{Row=435; F=83137; G=266},
    @{Row=436; F=139211; G=331},
    @{Row=438; F=118395; G=231},
    @{Row=439; F=86643; G=300},
    @{Row=440; F=72863; G=222},
    @{Row=441; F=65807; G=191},
    @{Row=442; F=67286; G=166},
    @{Row=443; F=102776; G=198},
    @{Row=444; F=100036; G=176},
    @{Row=446; F=86597; G=264},
    @{Row=447; F=64634; G=182},
    @{Row=448; F=58797; G=132},
    @{Row=449; F=59728; G=155},
    @{Row=450; F=87535; G=165},
    @{Row=451; F=82384; G=108},
    @{Row=453; F=67290; G=202},
    @{Row=454; F=50726; G=125},
    @{Row=455; F=49993; G=119},
    @{Row=456; F=47896; G=128},
    @{Row=457; F=75324; G=124},
    @{Row=458; F=67621; G=68},
    @{Row=460; F=55644; G=145},
    @{Row=461; F=43526; G=93},
    @{Row=462; F=41886; G=46},
    @{Row=463; F=44785; G=67},
    @{Row=464; F=69788; G=79},
    @{Row=465; F=58146; G=50},
    @{Row=467; F=50148; G=73},
    @{Row=468; F=40495; G=44},
    @{Row=469; F=39118; G=37},
    @{Row=470; F=41439; G=41},
    @{Row=471; F=62580; G=50},
    @{Row=472; F=47526; G=21},
    @{Row=473; F=38713; G=38},
    @{Row=474; F=43474; G=56},
    @{Row=475; F=33820; G=26},
    @{Row=476; F=34711; G=30},
    @{Row=477; F=36472; G=33},
    @{Row=478; F=49969; G=32},
    @{Row=479; F=38689; G=31},
    @{Row=480; F=31864; G=21},
    @{Row=481; F=41041; G=36},
    @{Row=482; F=31589; G=23}
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 6).Value = $u.F
    $ws.Cells.Item($u.Row, 7).Value = $u.G
}

# --- 3) Append the new row for 2021-06-30 ---
$newRow = 483
$ws.Range("A$newRow").Value = 44377
$ws.Range("B$newRow").Value = 391659
$ws.Range("C$newRow").Value = 5384
$ws.Range("D$newRow").Value = 17
$ws.Range("E$newRow").Value = 12511
$ws.Range("F$newRow").Value = 52016
$ws.Range("G$newRow").Value = 32
